$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert two new rows of market data right before row 420,
# pushing all subsequent rows (old 420-505) down to (422-507).
$ws.Rows("420:421").Insert()

# Row 420: new "Morada(o)" variety entry for 2022-03-21
$ws.Cells.Item(420, 1).Value = 5
$ws.Cells.Item(420, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(420, 3).Value = "Maule"
$ws.Cells.Item(420, 4).Value = 44641
$ws.Cells.Item(420, 5).Value = 7
$ws.Cells.Item(420, 6).Value = 100112004
$ws.Cells.Item(420, 7).Value = "Cebolla"
$ws.Cells.Item(420, 8).Value = "Morada(o)"
$ws.Cells.Item(420, 9).Value = "1a (guarda)"
$ws.Cells.Item(420, 10).Value = 500
$ws.Cells.Item(420, 11).Value = 8000
$ws.Cells.Item(420, 12).Value = 8000
$ws.Cells.Item(420, 13).Value = 8000
$ws.Cells.Item(420, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(420, 15).Value = "Región del Maule"
$ws.Cells.Item(420, 16).Value = 533
$ws.Cells.Item(420, 17).Value = 15
$ws.Cells.Item(420, 18).Value = "Hortaliza"

# Row 421: new "Sin especificar / 1a (cosecha)" entry for 2022-03-21
$ws.Cells.Item(421, 1).Value = 5
$ws.Cells.Item(421, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(421, 3).Value = "Maule"
$ws.Cells.Item(421, 4).Value = 44641
$ws.Cells.Item(421, 5).Value = 7
$ws.Cells.Item(421, 6).Value = 100112004
$ws.Cells.Item(421, 7).Value = "Cebolla"
$ws.Cells.Item(421, 8).Value = "Sin especificar"
$ws.Cells.Item(421, 9).Value = "1a (cosecha)"
$ws.Cells.Item(421, 10).Value = 2500
$ws.Cells.Item(421, 11).Value = 4500
$ws.Cells.Item(421, 12).Value = 4500
$ws.Cells.Item(421, 13).Value = 4500
$ws.Cells.Item(421, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(421, 15).Value = "Región del Maule"
$ws.Cells.Item(421, 16).Value = 180
$ws.Cells.Item(421, 17).Value = 25
$ws.Cells.Item(421, 18).Value = "Hortaliza"
